$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new date columns (Jun_15, Jun_17) immediately to the left of the
# existing "Jun_13"/"Jun_10" columns (B:C). This shifts the old B -> D and the
# old C -> E, carrying their values/formatting along automatically.
$ws.Columns("B:C").Insert()

# New header cells for the inserted columns.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Default rating value ("UN") for every data row in the two new columns,
# matching what the other date columns already contain.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
}

# Match the column widths of the other data columns (custom width, 8 chars).
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
